# Automatic update of files.
# Update the "Förändrad" date column (C) for rows 2-12 from 2023-09-03
# (serial 45172) to 2023-09-06 (serial 45175), keeping existing
# number formatting/style intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSerial = 45175

for ($row = 2; $row -le 12; $row++) {
    $ws.Range("C$row").Value = $newSerial
}
